# Refresh the cryptocurrency price/volume table (columns B:E, rows 2-51)
# on the active worksheet with the latest scraped data. The coin at the
# bottom of the previous ranking (Frax) dropped off and Decentraland was
# appended, shifting every subsequent row up by one position.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '30.449.80', '  +0.99%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.870.97', '  +0.55%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.000', '  +0.09%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '247.16', '  +2.41%  '),
    @(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  +0.06%  '),
    @(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.4744', '  +1.09%  '),
    @(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2907', '  +1.82%  '),
    @(9, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.06499', '  +0.43%  '),
    @(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '21.97', '  +6.50%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07719', '  +0.62%  '),
    @(12, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '97.71', '  +4.05%  '),
    @(13, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.7380', '  +8.46%  '),
    @(14, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.871.61', '  +0.04%  '),
    @(15, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '5.117', '  +1.23%  '),
    @(16, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '273.62', '  +1.29%  '),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '30.431.29', '  +0.98%  '),
    @(18, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '13.40', '  +0.97%  '),
    @(19, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.000007552', '  +0.24%  '),
    @(20, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '1.000', '  +0.07%  '),
    @(21, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.123.74', '  +1.29%  '),
    @(22, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.000', '  +0.11%  '),
    @(23, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '5.226', '  +1.42%  '),
    @(24, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '6.177', '  +1.62%  '),
    @(25, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '9.269', '  -0.65%  '),
    @(26, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '163.70', '  -1.74%  '),
    @(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '18.83', '  +0.66%  '),
    @(28, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.932', '  +2.76%  '),
    @(29, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1005', '  +2.19%  '),
    @(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.366', '  -0.69%  '),
    @(31, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.512', '  +1.07%  '),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.297', '  +1.97%  '),
    @(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.129', '  +3.32%  '),
    @(34, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.04817', '  +2.75%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.121', '  +1.06%  '),
    @(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6979', '  +2.04%  '),
    @(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.709', '  -0.21%  '),
    @(38, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01860', '  +2.23%  '),
    @(39, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.749', '  +0.99%  '),
    @(40, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.264', '  -1.28%  '),
    @(41, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.972', '  +4.88%  '),
    @(42, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '71.15', '  +1.88%  '),
    @(43, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4196', '  +3.68%  '),
    @(44, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9999', '  +0.07%  '),
    @(45, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.8363', '  +0.17%  '),
    @(46, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '102.93', '  +0.97%  '),
    @(47, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.326', '  +1.42%  '),
    @(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.021', '  +1.43%  '),
    @(49, 'Elrond', 'https://coinranking.com/coin/omwkOTglq+elrond-egld', '35.45', '  +3.57%  '),
    @(50, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '922.27', '  -0.34%  '),
    @(51, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '0.3884', '  +3.41%  ')
)

foreach ($row in $data) {
    $r = $row[0]

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Column D ("Price") can look like a plain number (e.g. "1.000" or
    # "247.16"), which Excel would otherwise silently reinterpret as a
    # numeric value and reformat/round. Force it to be stored as text,
    # matching the original inline-string cells, then drop the
    # temporary "Text" number format so no extra cell style lingers.
    $priceCell = $ws.Cells.Item($r, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $row[3]
    $priceCell.ClearFormats()

    $ws.Cells.Item($r, 5).Value = $row[4]
}
